{"js": "// Apply the \"Plan\" document edits described by the commit:\n//  1. Fix typo/count: \"Start with 3 service\" -> \"Start with 2 services\"\n//  2. Highlight two risk-call-out bullets in red font color\n//  3. Add a new \"risk management\" checklist section at the end of the\n//     \"Checklist / not sure where to put them\" list\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1. Fix the bullet text \"Start with 3 service\" -> \"Start with 2 services\"\n// ---------------------------------------------------------------------\nconst startWithResults = body.search(\"Start with 3 service\", { matchCase: true });\nstartWithResults.load(\"text\");\nawait context.sync();\n\nif (startWithResults.items.length > 0) {\n  startWithResults.items[0].insertText(\"Start with 2 services\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2. Colour the two bullets about re-try logic / end-to-end testing red\n// ---------------------------------------------------------------------\nconst redTexts = [\n  \"Don\\u2019t worry about re-try logic\",\n  \"End-to-end test with two running remote instances would be nice\"\n];\n\nfor (const redText of redTexts) {\n  const results = body.search(redText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].font.color = \"#FF0000\";\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3. Insert new \"risk management\" checklist items after\n//    \"Logging certain events ...\" bullet (still inside the same\n//    List Paragraph / numId 5 bulleted list).\n// ---------------------------------------------------------------------\nconst anchorResults = body.search(\n  \"Logging certain events\",\n  { matchCase: true }\n);\nanchorResults.load(\"text\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  let anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n\n  const newItems = [\n    { text: \"Have I applied risk management?\", level: 0 },\n    { text: \"\\u201CWhat if this goes wrong?\\u201D\", level: 1 },\n    { text: \"\\u201CWhat if this hardware/software component stops working?\\u201D\", level: 1 },\n    { text: \"\\u201CWhat if a queue goes down?\\u201D\", level: 1 }\n  ];\n\n  for (const item of newItems) {\n    const newParagraph = anchorParagraph.insertParagraph(item.text, \"After\");\n    newParagraph.listItemOrNullObject.level = item.level;\n    await context.sync();\n    anchorParagraph = newParagraph;\n  }\n}\n", "ps1": "# Apply the \"Plan\" document edits described by the commit:\n#  1. Fix typo/count: \"Start with 3 service\" -> \"Start with 2 services\"\n#  2. Highlight two risk-call-out bullets in red font color\n#  3. Add a new \"risk management\" checklist section at the end of the\n#     \"Checklist / not sure where to put them\" list\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. Fix the bullet text \"Start with 3 service\" -> \"Start with 2 services\"\n#    (Assign to the paragraph Range.Text directly - rather than using\n#    Find/Replace - so the run's empty <w:rPr/> element is preserved.)\n# ---------------------------------------------------------------------\n$startWithParagraphs = $d.Paragraphs\nfor ($i = 1; $i -le $startWithParagraphs.Count; $i++) {\n    if ($startWithParagraphs.Item($i).Range.Text -like \"*Start with 3 service*\") {\n        $startWithParagraphs.Item($i).Range.Text = \"Start with 2 services\"\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 2. Colour the two bullets about re-try logic / end-to-end testing red\n# ---------------------------------------------------------------------\n$redTexts = @(\n    \"Don't worry about re-try logic\",\n    \"End-to-end test with two running remote instances would be nice\"\n)\n\nforeach ($redText in $redTexts) {\n    $range = $d.Content\n    $rfind = $range.Find\n    $rfind.ClearFormatting()\n    $rfind.Text = $redText\n    $found = $rfind.Execute($rfind.Text, $false, $true, $false, $false, $false, $true, 1, $false)\n    if ($found) {\n        $range.Font.Color = 255\n    }\n}\n\n# ---------------------------------------------------------------------\n# 3. Insert new \"risk management\" checklist items after\n#    \"Logging certain events ...\" bullet (still inside the same\n#    List Paragraph / numId 5 bulleted list). We locate the anchor\n#    paragraph by scanning $d.Paragraphs because Find.Execute collapses\n#    its Range to the matched text only, which is not useful for\n#    retrieving the enclosing paragraph via .Paragraphs.Item(1).\n# ---------------------------------------------------------------------\n$allParagraphs = $d.Paragraphs\n$anchorIndex = -1\nfor ($i = 1; $i -le $allParagraphs.Count; $i++) {\n    if ($allParagraphs.Item($i).Range.Text -like \"*Logging certain events*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -gt 0) {\n    $quote = [char]0x201C\n    $endQuote = [char]0x201D\n\n    $newItemTexts = @(\n        \"Have I applied risk management?\",\n        ($quote + \"What if this goes wrong?\" + $endQuote),\n        ($quote + \"What if this hardware/software component stops working?\" + $endQuote),\n        ($quote + \"What if a queue goes down?\" + $endQuote)\n    )\n    # ListLevelNumber is 1-based: 1 = ilvl 0 (top item), 2 = ilvl 1 (sub item)\n    $newItemLevels = @(1, 2, 2, 2)\n\n    $currentAnchorParagraph = $allParagraphs.Item($anchorIndex)\n\n    for ($i = 0; $i -lt $newItemTexts.Length; $i++) {\n        $currentAnchorParagraph.Range.InsertParagraphAfter()\n        $anchorIndex = $anchorIndex + 1\n        $refreshedParagraphs = $d.Paragraphs\n        $currentAnchorParagraph = $refreshedParagraphs.Item($anchorIndex)\n        $currentAnchorParagraph.Range.Text = $newItemTexts[$i]\n        $currentAnchorParagraph.Range.ListFormat.ListLevelNumber = $newItemLevels[$i]\n    }\n}\n"}
